$d = $word.ActiveDocument

$pairs = @(
    @("632÷4=158, 0", "124÷5=24, 4"),
    @("475÷5=95, 0", "759÷8=94, 7"),
    @("639÷8=79, 7", "787÷8=98, 3"),
    @("886÷7=126, 4", "887÷7=126, 5"),
    @("265÷9=29, 4", "101÷7=14, 3"),
    @("845÷7=120, 5", "873÷6=145, 3"),
    @("304÷9=33, 7", "886÷3=295, 1"),
    @("780÷3=260, 0", "769÷2=384, 1"),
    @("665÷6=110, 5", "924÷6=154, 0"),
    @("489÷9=54, 3", "299÷7=42, 5"),
    @("533÷3=177, 2", "407÷6=67, 5"),
    @("634÷7=90, 4", "930÷9=103, 3"),
    @("153÷3=51, 0", "293÷2=146, 1"),
    @("555÷7=79, 2", "134÷6=22, 2"),
    @("707÷7=101, 0", "445÷2=222, 1"),
    @("512÷3=170, 2", "724÷9=80, 4"),
    @("692÷8=86, 4", "772÷9=85, 7"),
    @("184÷2=92, 0", "728÷8=91, 0"),
    @("521÷7=74, 3", "910÷9=101, 1"),
    @("564÷3=188, 0", "988÷7=141, 1"),
    @("842÷5=168, 2", "708÷4=177, 0"),
    @("820÷5=164, 0", "431÷4=107, 3"),
    @("172÷7=24, 4", "432÷9=48, 0"),
    @("584÷8=73, 0", "472÷2=236, 0"),
    @("990÷7=141, 3", "403÷6=67, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
